$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ryddet litt i seminaroppgavene: hvert oppgaveseminar fikk et lopenummer
# (1-7), og "hullet" i kalenderen (rad 8 / D8, som tidligere var
# "Ingen aktivitet i Auditorium") er fylt med seminar 4 den 06.10.
$ws.Range("D3").Value  = "01.09: Oppgaveseminar 1 Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D5").Value  = "15.09: Oppgaveseminar 2 Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D7").Value  = "22.09: Oppgaveseminar 3 Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D8").Value  = "06.10: Oppgaveseminar 4 Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D10").Value = " 20.10: Oppgaveseminar 5 Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D12").Value = " 03.11: Oppgaveseminar  6, Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D13").Value = "10.11: Oppgaveseminar 7, Aud A. Se \@ref(seminar) for oppgaver."

# Flytt den aktive markoren, slik den sto etter redigeringen.
[void]$ws.Range("C22").Select()
